$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.723.72"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.759.56"
$ws.Range("E3").Value = "  -1.98%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.82"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4421"
$ws.Range("E7").Value = "  -2.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3739"
$ws.Range("E8").Value = "  +0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.42"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07703"
$ws.Range("E10").Value = "  +1.95%  "

$ws.Range("E11").Value = "  -1.44%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.72"
$ws.Range("E13").Value = "  -2.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.191"
$ws.Range("E14").Value = "  -1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.419"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.755.81"
$ws.Range("E16").Value = "  -1.92%  "

$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.58"
$ws.Range("E18").Value = "  +10.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06225"
$ws.Range("E19").Value = "  -7.85%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.34"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5320"
$ws.Range("E23").Value = "  -3.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.753.53"
$ws.Range("E24").Value = "  -1.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.64"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.314"
$ws.Range("E26").Value = "  -4.13%  "

$ws.Range("E27").Value = "  +1.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.54"
$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.956.16"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.46"
$ws.Range("E31").Value = "  -3.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.213"
$ws.Range("E32").Value = "  -1.90%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.761"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09302"
$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.641"
$ws.Range("E35").Value = "  -9.85%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.68"
$ws.Range("E36").Value = "  +5.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2189"
$ws.Range("E37").Value = "  -7.65%  "

$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06147"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6498"
$ws.Range("E40").Value = "  -0.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.084"
$ws.Range("E41").Value = "  -2.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.200"
$ws.Range("E42").Value = "  -0.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.995"
$ws.Range("E43").Value = "  -4.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.420"
$ws.Range("E44").Value = "  -4.02%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.80"
$ws.Range("E46").Value = "  -2.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6017"
$ws.Range("E47").Value = "  -0.91%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.767"
$ws.Range("E48").Value = "  -1.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.38"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.999"
$ws.Range("E50").Value = "  -1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.144"
$ws.Range("E51").Value = "  -1.17%  "
